# Rename the worksheet from "adactin" to "com.adactin.hotelapp".
# Excel automatically keeps the _xlnm._FilterDatabase defined name (and any
# other formulas) in sync with the new sheet name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "com.adactin.hotelapp"
